$d = $word.ActiveDocument

# The document starts with a "title block" paragraph:
#   Title: 3YCM_report
#   Author: Sarah Stevens
# This whole paragraph (including its paragraph mark) needs to be removed,
# so that the document now begins directly with the "Third Year Committee
# Report" heading.

$titlePara = $d.Paragraphs(1)
$titlePara.Range.Delete()
